# CW3M 1.2.1 release: add No_wetlands_demo scenario rows to the
# "2010 and 2010-18" regression-testing sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows after the current row 111 (the existing blank
# separator row stays put at 111, absorbing the new "CW3M C650" data row;
# two more rows open up under it for the "No_wetlands_demo" data row and a
# fresh pair of blank separator rows). This pushes the old rows 113-115
# down to 116-118 and removes the old row-112 gap.
$ws.Rows("112:114").Insert()

# ---- Row 111: new CW3M C650 baseline data row ---------------------------
$ws.Cells.Item(111, 1).Value2 = "CW3M C650"
$ws.Cells.Item(111, 2).Value2 = "Demo_Baseline_2010-18_C650"
$ws.Cells.Item(111, 3).Value2 = "2010-18"
$ws.Cells.Item(111, 4).Value2 = 936.4586451111112
$ws.Cells.Item(111, 5).Value2 = 1890.2624918888889
$ws.Cells.Item(111, 6).Value2 = 1.1353011111111111
$ws.Cells.Item(111, 7).Value2 = 270.41205844444437
$ws.Cells.Item(111, 8).Value2 = 9.8446185555555559
$ws.Cells.Item(111, 9).Value2 = 7.3262146666666679
$ws.Cells.Item(111, 10).Value2 = 8.202840444444444
$ws.Cells.Item(111, 11).Value2 = 664.89181855555546
$ws.Cells.Item(111, 12).Value2 = 80.365177222222229
$ws.Cells.Item(111, 13).Value2 = 1422.6016167777777
$ws.Cells.Item(111, 14).Value2 = 939.5350204444444
$ws.Cells.Item(111, 15).Value2 = 6416.2319064444446
$ws.Cells.Item(111, 16).Value2 = 27412.947482666666
$ws.Cells.Item(111, 17).Value2 = 0.15714333333333327
$ws.Cells.Item(111, 18).Value2 = [double]"4.0333333333333261E-5"

$ws.Range("D111:N111").NumberFormat = "0.00"
$ws.Range("O111:P111").NumberFormat = "0"
$ws.Range("Q111").NumberFormat = "0.00"
$ws.Range("R111").NumberFormat = "0.000000"

# ---- Row 112: new No_wetlands_demo data row ------------------------------
$ws.Cells.Item(112, 1).Value2 = "CW3M C650"
$ws.Cells.Item(112, 2).Value2 = "No_wetlands_demo_2010-18_C650"
$ws.Cells.Item(112, 3).Value2 = "2010-18"
$ws.Cells.Item(112, 4).Value2 = 934.65285922222233
$ws.Cells.Item(112, 5).Value2 = 1890.2624918888889
$ws.Cells.Item(112, 6).Value2 = 1.0305982222222221
$ws.Cells.Item(112, 7).Value2 = 270.41205844444437
$ws.Cells.Item(112, 8).Value2 = 9.860652222222221
$ws.Cells.Item(112, 9).Value2 = 7.342017444444445
$ws.Cells.Item(112, 10).Value2 = 8.2161996666666663
$ws.Cells.Item(112, 11).Value2 = 664.22806799999989
$ws.Cells.Item(112, 12).Value2 = 80.302231777777777
$ws.Cells.Item(112, 13).Value2 = 1423.3923340000001
$ws.Cells.Item(112, 14).Value2 = 937.57914911111118
$ws.Cells.Item(112, 15).Value2 = 5596.6441785555553
$ws.Cells.Item(112, 16).Value2 = 27457.329643999998
$ws.Cells.Item(112, 17).Value2 = 0.15730522222222224
$ws.Cells.Item(112, 18).Value2 = [double]"2.9444444444444438E-5"

$ws.Range("D112:N112").NumberFormat = "0.00"
$ws.Range("O112:P112").NumberFormat = "0"
$ws.Range("Q112").NumberFormat = "0.00"
$ws.Range("R112").NumberFormat = "0.000000"
$ws.Range("O112:P112").Interior.Color = 65535

# ---- Row 113/114: blank separator rows (highlight O:P like row 112) -----
$ws.Range("O113:P114").Interior.Color = 65535

# ---- Selection / view state to match the saved workbook ------------------
$ws.Range("A113:XFD114").Select()
